$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 12 / column G: change full-width colons to half-width colons
$ws.Range("G12").Value = "0:開帳`n1:關帳`n2:關帳取消"

# Row 16 / column G: change colon to equals sign in the parenthetical codes
$ws.Range("G16").Value = "只更新特定筆(09:放款)`n預設為000，產生上傳媒體(02=支票繳款，09=放款)關帳時＋１"

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("H16").Select()
